$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing data
# (Player ID ... etc.) one column to the right.
$ws.Columns("A").Insert()

# New header cell for the inserted "Match ID" column (row 2 is the header row).
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Row 3 is a hidden spacer row; give the new A3 cell the same bold style
# without disturbing the stored row height (temporarily unhide while writing).
$ws.Rows(3).Hidden = $false
$ws.Range("A3").Font.Bold = $true
$ws.Rows(3).Hidden = $true

# Data rows 4-14 all belong to match id 15.
$ws.Range("A4:A14").Value = 15
$ws.Range("A4:A14").Font.Bold = $true

# Row 15 is the hidden totals row; unhide while writing to avoid a spurious
# row-height change, then re-hide.
$ws.Rows(15).Hidden = $false
$ws.Range("A15").Value = 15
$ws.Rows(15).Hidden = $true

# Restore the visible selection to the newly added Match ID column.
[void]$ws.Range("A2:A14").Select()
